# Insert a new data row at row 13 (pushing the existing rows 13-66 down to 14-67),
# then populate the new row 13 as a duplicate of the (now shifted) original row 13
# data with updated Fecha (D), Variedad (K) and Origen (R) values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row before row 13; existing row 13 (and below) shift to row 14.
$ws.Rows.Item(13).Insert()

# Copy the now-shifted original row (row 14) into the new blank row 13.
$ws.Rows.Item(14).Copy()
$ws.Rows.Item(13).PasteSpecial()
$excel.CutCopyMode = $false

# Overwrite the fields that differ for the newly inserted record.
$ws.Cells.Item(13, 4).Value2 = 44453
$ws.Cells.Item(13, 11).Value = "Murcott"
$ws.Cells.Item(13, 18).Value = "Región de Coquimbo"
